# "Reports file name changed"
#
# The sheet used to track 5 test-run rows (Test Case ID, Test Case, Links,
# Start Time, End Time, Result). The report now only keeps the 3-column
# header row, and the "Links" header was renamed to "Pges".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the 3rd header ("Links" -> "Pges").
$ws.Range("C1").Value = "Pges"

# Drop all the data rows (2-6); only the header row remains.
$ws.Range("A2:F6").EntireRow.Delete()

# Drop the now-unused trailing headers (Start Time / End Time / Result)
# that lived in D1:F1 - use Clear() (not ClearContents) so the cells are
# fully removed from the sheet (and the used range shrinks to A1:C1)
# instead of being left behind as empty-but-styled cells.
$ws.Range("D1:F1").Clear()

# The remaining header row no longer needs its own explicit row-level
# format override - reproduce the look (bold font on a yellow fill) on
# just the 3 header cells that are left.
$ws.Rows(1).ClearFormats()
$ws.Range("A1:C1").Interior.Color = 65535
$ws.Range("A1:C1").Font.Bold = $true

# Park the selection like the saved file does.
[void]$ws.Range("D7").Select()
